# Apply the "rdbms" commit changes to the Nexial function catalog sheet ("#system").
#  1. Remove "clearVariables(variables)" from the "base" list (column F) - shrinks base from F2:F40 to F2:F39
#  2. Add "terminate(programName)" to the "external" list (column J) - grows external from J2:J5 to J2:J6
#  3. Add "saveSelectedText(var,locator)" and "saveSelectedValue(var,locator)" to the "web" list (column Z)
#     - grows web from Z2:Z135 to Z2:Z137
#
# NOTE: Range.Delete()/Insert() in this runtime shift the *entire row* rather than
# being confined to the targeted column, so the column contents are shifted manually,
# cell-by-cell, to avoid disturbing the other function-catalog columns on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1. Remove F19 ("clearVariables(variables)") by shifting F19:F40 up by one cell ---
for ($r = 19; $r -le 39; $r++) {
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r + 1, 6).Value2
}
$ws.Cells.Item(40, 6).ClearContents()

# --- 2. Add new "external" entry at J6 (first empty row right after current J2:J5 list) ---
$ws.Cells.Item(6, 10).Value2 = "terminate(programName)"

# --- 3. Make room for two new rows in column Z only, by shifting Z99:Z135 down to Z101:Z137 ---
for ($r = 135; $r -ge 99; $r--) {
    $ws.Cells.Item($r + 2, 26).Value2 = $ws.Cells.Item($r, 26).Value2
}
$ws.Cells.Item(99, 26).Value2 = "saveSelectedText(var,locator)"
$ws.Cells.Item(100, 26).Value2 = "saveSelectedValue(var,locator)"

# --- 4. Update the workbook-level defined names to reflect the new ranges ---
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$39"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$6"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$137"
